# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 5
$ws1.Range("F5").Value = 37
$ws1.Range("F7").Value = 1787
$ws1.Range("F11").Value = 2109
$ws1.Range("F12").Value = 36
$ws1.Range("F13").Value = 147
$ws1.Range("F14").Value = 1349
$ws1.Range("F15").Value = 473
$ws1.Range("F18").Value = 211
$ws1.Range("F24").Value = 16
$ws1.Range("F25").Value = 1140
$ws1.Range("F27").Value = 340
$ws1.Range("F29").Value = 274
$ws1.Range("F30").Value = 335

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5
$ws4.Range("F5").Value = 37
$ws4.Range("F7").Value = 1787
$ws4.Range("F12").Value = 2109
$ws4.Range("F13").Value = 36
$ws4.Range("F14").Value = 147
$ws4.Range("F15").Value = 1349
$ws4.Range("F16").Value = 473
$ws4.Range("F19").Value = 211
$ws4.Range("F25").Value = 16
$ws4.Range("F26").Value = 1140
$ws4.Range("F28").Value = 340
$ws4.Range("F30").Value = 274
$ws4.Range("F31").Value = 335
